$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.032.99'
$ws.Range("E2").Value = '  +0.23%  '

$ws.Range("D3").Value = '1.922.01'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.33%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.37'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4597'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3814'
$ws.Range("E8").Value = '  -0.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07752'
$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9778'
$ws.Range("E10").Value = '  +0.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.74'
$ws.Range("E11").Value = '  +2.49%  '

$ws.Range("D12").Value = '1.936.17'
$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.706'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.965'
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07007'
$ws.Range("E15").Value = '  -0.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.40'
$ws.Range("E17").Value = '  +0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009523'
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.71'
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").Value = '29.021.13'
$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.353'
$ws.Range("E22").Value = '  +0.27%  '

$ws.Range("E23").Value = '  +0.72%  '

$ws.Range("D24").Value = '2.142.57'
$ws.Range("E24").Value = '  +0.83%  '

$ws.Range("E25").Value = '  +0.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.46'
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.07'
$ws.Range("E27").Value = '  -0.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.650'
$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.99'
$ws.Range("E29").Value = '  +0.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.836'
$ws.Range("E30").Value = '  +0.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09329'
$ws.Range("E31").Value = '  +0.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8613'
$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.098'
$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.238'
$ws.Range("E34").Value = '  -0.96%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.018'
$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05694'
$ws.Range("E36").Value = '  -0.44%  '

$ws.Range("E37").Value = '  +0.93%  '

$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02044'
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.097'
$ws.Range("E40").Value = '  +14.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.432'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5487'
$ws.Range("E42").Value = '  -0.64%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1755'
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.387'
$ws.Range("E44").Value = '  +1.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002827'
$ws.Range("E45").Value = '  -0.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.189'
$ws.Range("E46").Value = '  +5.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5190'
$ws.Range("E47").Value = '  -0.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.24'
$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06924'
$ws.Range("E49").Value = '  +1.56%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.38'
$ws.Range("E50").Value = '  -1.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.762'
$ws.Range("E51").Value = '  -1.32%  '
